$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row labels: "_old" suffix -> "_FV2410", "_new" suffix -> "_FV2504"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2410"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2504"
        }
    }
}

# 2. Convert the data range into an Excel Table (ListObject)
$range = $ws.Range("A1:U81")
$listObject = $ws.ListObjects.Add(1, $range, $null, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# 3. Freeze the header row (pane split after row 1)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
